$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the long note in A20: "VC11 - batch build..." -> "VC11 / VC12 - batch build..."
$ws.Range("A20").Value2 = "VC11 / VC12 - batch build / rebuild of QuantLibAddin - kick it off in the evening - following morning it's still running."

# Extend the "TO" (Timed Out) marker from J6:Q6 to also cover R6:Y6 (new VC12 columns),
# matching the style used by the existing TO cells.
$ws.Range("R6:Y6").Value2 = "TO"
$ws.Range("J6").Copy()
$ws.Range("R6:Y6").PasteSpecial(-4122)
$excel.CutCopyMode = 0
